$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "机构号" (institution number) column for the two data rows:
# previously "上海FTU_GM1", now "GM_151"
$ws.Range("X2").Value = "GM_151"
$ws.Range("X3").Value = "GM_151"
